# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-column suffixes to "_FV2410" / "_FV2504"
# - Turn the header + data range into an Excel Table ("Table1")
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels (row 1) ---------------------------------
# Columns A:J were "<Label>_old" -> "<Label>_FV2410"
# Column K ("diff") is unchanged
# Columns L:U were "<Label>_new" -> "<Label>_FV2504"

$labels = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $labels.Count; $i++) {
    $oldCol = [char](65 + $i)        # A..J
    $newCol = [char](76 + $i)        # L..U
    $ws.Range("$oldCol" + "1").Value = $labels[$i] + "_FV2410"
    $ws.Range("$newCol" + "1").Value = $labels[$i] + "_FV2504"
}

# --- 2. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the A1:U66 range into an Excel Table -----------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), $null, 1)
$tbl.Name = "Table1"

Write-Host "Header labels renamed, freeze panes applied, Table1 created."
